$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated 2D training schedule data (rows 2-6), no break screen (no I=4 "break" values; now I=5 throughout)
$data = @(
    @(1, 9, 2, 6, 5, -3, 3, 43, 5),
    @(2, 5, 0, 0, 1, -5, 1, 65, 5),
    @(3, 8, 1, 7, 6, -1, 5, 21, 5),
    @(4, 5, 2, 1, 4, -4, 2, 54, 5),
    @(5, 9, 1, 7, 5, -2, 4, 32, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

$ws.Range("I1").Select() | Out-Null
